$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Resize/reposition the PivotChart's graphic frame so the two-cell anchor
# changes from (from col2/27214,row0/40821)-(to col13/0,row16/15875)
# to (from col2/27215,row0/40821)-(to col12/559595,row15/345281).
$co = $ws.ChartObjects().Item(1)
$co.Left = 188.12240557332677
$co.Top = 3.214251968503937
$co.Width = 626.29468503937
$co.Height = 473.97322834645666

# Change the zoom level of the active sheet view from 60% to 80%.
$ws.Activate()
$excel.ActiveWindow.Zoom = 80
